# Rename CAPACITY_COST_* / DISPATCH_COST_* (and related) keywords to
# FIXED_COST_* / VAR_COST_* throughout every worksheet of the workbook.
#
# This mirrors the commit "changed CAPACITY and DISPATCH to FIXED and VAR
# in case_input.*" -- a pure text rename of keyword tokens (and the two
# explanatory sentences that mention them) that appear in column A (and a
# few header cells further to the right) of each of the four sheets.

$wb = $excel.ActiveWorkbook

# Exact-match keyword renames (whole-cell-value replacements).
$exactMap = @{
    "SOLAR_CAPACITY_FILE"             = "SOLAR_FIXED_FILE"
    "CAPACITY_COST_SOLAR"             = "FIXED_COST_SOLAR"
    "DISPATCH_COST_SOLAR"             = "VAR_COST_SOLAR"
    "WIND_CAPACITY_FILE"              = "WIND_FIXED_FILE"
    "CAPACITY_COST_WIND"              = "FIXED_COST_WIND"
    "DISPATCH_COST_WIND"              = "VAR_COST_WIND"
    "CAPACITY_COST_NATGAS"            = "FIXED_COST_NATGAS"
    "DISPATCH_COST_NATGAS"            = "VAR_COST_NATGAS"
    "CAPACITY_COST_NUCLEAR"           = "FIXED_COST_NUCLEAR"
    "DISPATCH_COST_NUCLEAR"           = "VAR_COST_NUCLEAR"
    "CAPACITY_COST_STORAGE"           = "FIXED_COST_STORAGE"
    "DISPATCH_COST_TO_STORAGE"        = "VAR_COST_TO_STORAGE"
    "DISPATCH_COST_FROM_STORAGE"      = "VAR_COST_FROM_STORAGE"
    "DISPATCH_COST_UNMET_DEMAND"      = "VAR_COST_UNMET_DEMAND"
    "CAPACITY_COST_PGP_STORAGE"       = "FIXED_COST_PGP_STORAGE"
    "DISPATCH_COST_TO_PGP_STORAGE"    = "VAR_COST_TO_PGP_STORAGE"
    "DISPATCH_COST_FROM_PGP_STORAGE"  = "VAR_COST_FROM_PGP_STORAGE"
    "CAPACITY_COST_TO_PGP_STORAGE"    = "FIXED_COST_TO_PGP_STORAGE"
    "CAPACITY_COST_FROM_PGP_STORAGE"  = "FIXED_COST_FROM_PGP_STORAGE"
}

# Longer explanatory sentences that reference the old lower-case keyword
# names inline; handled via substring replacement instead of exact match.
$substringFrom = @(
    "capacity_cost_wind",
    "dispatch_cost_wind",
    "dispatch_cost_storage"
)
$substringTo = @(
    "FIXED_cost_wind",
    "VAR_cost_wind",
    "VAR_cost_storage"
)

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
            $val = $cell.Value2

            if ($null -eq $val) { continue }
            if ($val.GetType().Name -ne "String") { continue }

            if ($exactMap.ContainsKey($val)) {
                $cell.Value = $exactMap[$val]
                continue
            }

            $newVal = $val
            for ($i = 0; $i -lt $substringFrom.Count; $i++) {
                $newVal = $newVal.Replace($substringFrom[$i], $substringTo[$i])
            }
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
